$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove existing hyperlinks up front - they do not travel with cells when
#    columns are inserted, so we recreate all of them (old + new) afterwards
#    pointing at their final locations.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2. Insert two new columns: "UserName" (new B) and "Confirm Password" (new D)
#    Doing this in sequence shifts every existing column from old-B onward by
#    two positions (old B -> new C, old C -> new E, old D -> new F, etc.)
# ---------------------------------------------------------------------------
$ws.Range("B1").EntireColumn.Insert()
$ws.Range("D1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 3. New column headers / values
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "UserName"
$ws.Range("D1").Value = "Confirm Password"

$ws.Range("B2").Value = "qatesting.lotuswave@gmail.com"
$ws.Range("D2").Value = "Lotuswave@123"

# ---------------------------------------------------------------------------
# 4. New appended columns at the end (U / V)
# ---------------------------------------------------------------------------
$ws.Range("U1").Value = "Quantity"
$ws.Range("V1").Value = "Discountcode"

$ws.Range("U6").Value = 3
$ws.Range("V7").Value = "GGQA`$25"

# ---------------------------------------------------------------------------
# 5. New row 5 - "Invalid details" login-style block
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "Invalid details"
$ws.Range("C5").Value = "Lotuswave"
$ws.Range("D5").Value = "Lotus123"
$ws.Range("E5").Value = "Testing"
$ws.Range("F5").Value = "Test"
$ws.Range("G5").Value = "qaflask.com"

# ---------------------------------------------------------------------------
# 6. New row 6 - "Product Qunatity" label (H6 blank cell already shifted from
#    the old F6 via the column insert above)
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "Product Qunatity"

# ---------------------------------------------------------------------------
# 7. New row 7 - "Discount" label + new D7 blank cell
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "Discount"
$ws.Range("D7").Value = ""

# ---------------------------------------------------------------------------
# 8. New rows 8-10 - credit card test data
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "CCMastercard"
$ws.Range("R8").Value = "5555555555554444"
$ws.Range("S8").Value = "04/26"
$ws.Range("T8").Value = 123

$ws.Range("A9").Value = "CCAmexcard"
$ws.Range("R9").Value = "378282246310005"
$ws.Range("S9").Value = "04/26"
$ws.Range("T9").Value = 1234

$ws.Range("A10").Value = "CCDiscovercard"
$ws.Range("R10").Value = "6011111111111117"
$ws.Range("S10").Value = "04/26"
$ws.Range("T10").Value = 123

# ---------------------------------------------------------------------------
# 9. Row 15 - new D15 blank cell (C15 already shifted from the old B15)
# ---------------------------------------------------------------------------
$ws.Range("D15").Value = ""

# ---------------------------------------------------------------------------
# 10. Recreate hyperlinks at their final resting places.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:Lotuswave@123", $null, $null, $null)
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:qatesting.lotuswave@gmail.com", $null, $null, $null)
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:Lotuswave@123", $null, "Lotuswave@123", $null)
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:qatesting.lotuswave@gmail.com", $null, $null, $null)
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:Lotuswave@123", $null, $null, $null)

# ---------------------------------------------------------------------------
# 11. Sheet view / window cosmetics from the diff
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 90
$ws.Range("O14").Select()

Write-Host "done"
